$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OTIS")

# Row 6: "Change in inventories"
$ws.Range("B6").Value = -66000000.0
$ws.Range("C6").Value = -76000000.0
$ws.Range("D6").Value = -72000000.0
$ws.Range("E6").Value = -35000000.0
$ws.Range("F6").Value = 4000000.0

# Row 8: "Change in payables and accrued liability"
$ws.Range("B8").Value = 4672000000.0
$ws.Range("C8").Value = 2924000000.0
$ws.Range("D8").Value = 1531000000.0
$ws.Range("E8").Value = 201000000.0
$ws.Range("F8").Value = -14000000.0
